# Tokenizer coverage refactor:
#  - paragraphs 10, 11, 14 (1-based Paragraphs index): collapse the
#    multi-run / lang-tagged runs into a single plain run (no rPr) and
#    drop the paragraph-mark rPr (pPr) entirely.
#  - paragraphs 12, 13, 15: collapse the empty-but-formatted paragraphs
#    down to fully empty <w:p/> (no pPr/rPr at all).
#  - final "Douglas ..." paragraph: split the second run's text into a
#    run-per-word (plus separate double/triple-space runs) sequence.

function New-OpenXmlPackage([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$d = $word.ActiveDocument

$dq = [string]([char]34)
$lsq = [string]([char]0x2018)
$rsq = [string]([char]0x2019)

# Paragraph 10: "Quote Problem" (straight double quotes), 3 runs -> 1 run, no pPr.
$p10Xml = New-OpenXmlPackage ('<w:body><w:p><w:r><w:t>' + $dq + 'Quote Problem' + $dq + '</w:t></w:r></w:p></w:body>')
$d.Paragraphs(10).Range.InsertXML($p10Xml)

# Paragraph 11: 'Quote Problem' (curly single quotes), 3 runs -> 1 run, no pPr.
$p11Xml = New-OpenXmlPackage ('<w:body><w:p><w:r><w:t>' + $lsq + 'Quote Problem' + $rsq + '</w:t></w:r></w:p></w:body>')
$d.Paragraphs(11).Range.InsertXML($p11Xml)

# Paragraphs 12 & 13: empty formatted paragraphs -> fully empty <w:p/>.
$emptyXml = New-OpenXmlPackage '<w:body><w:p/></w:body>'
$d.Paragraphs(12).Range.InsertXML($emptyXml)
$d.Paragraphs(13).Range.InsertXML($emptyXml)

# Paragraph 14: [Bracket Problem] -> single run, no pPr.
$p14Xml = New-OpenXmlPackage '<w:body><w:p><w:r><w:t>[Bracket Problem]</w:t></w:r></w:p></w:body>'
$d.Paragraphs(14).Range.InsertXML($p14Xml)

# Paragraph 15: empty formatted paragraph -> fully empty <w:p/>.
$d.Paragraphs(15).Range.InsertXML($emptyXml)

# Final paragraph ("Douglas ... French at all."): split the second run
# (everything after "Douglas") into a run per word plus dedicated
# double/triple-space runs.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$full = $lastPara.Range
$douglasLen = "Douglas".Length
$r2 = $d.Range($full.Start + $douglasLen, $full.End - 1)

$words = @(' ', '  ', 'lives ', '  ', 'in ', '  ', 'Florida ', '  ', 'and ', '  ', 'speaks ', '  ', 'no ', '  ', 'French ', '  ', 'at ', '   ', 'all.')
$runsXml = ''
foreach ($w in $words) {
    if ($w -eq 'all.') {
        $runsXml += '<w:r><w:t>' + $w + '</w:t></w:r>'
    } else {
        $runsXml += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
    }
}
$tailXml = New-OpenXmlPackage ('<w:body><w:p>' + $runsXml + '</w:p></w:body>')
$r2.InsertXML($tailXml)
